$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40657
$ws.Range("J3").Value = 40657
$ws.Range("L3").Value = 40657
$ws.Range("N3").Value = -40885
$ws.Range("H41").Value = 695.6111
$ws.Range("J41").Value = 733.6667
$ws.Range("L41").Value = 733.6667
$ws.Range("N41").Value = -1613.6667
$ws.Range("H55").Value = 227.21053
$ws.Range("I55").Value = 108.57143
$ws.Range("J55").Value = 296.41666
$ws.Range("K55").Value = 108.57143
$ws.Range("L55").Value = 296.41666
$ws.Range("M55").Value = 105.42857
$ws.Range("N55").Value = -724.41666
$ws.Range("H74").Value = 8799.4
$ws.Range("I74").Value = 4749.25
$ws.Range("K74").Value = 4749.25
$ws.Range("M74").Value = -3813.25
$ws.Range("H77").Value = 8799.4
$ws.Range("I77").Value = 4749.25
$ws.Range("K77").Value = 23746.25
$ws.Range("M77").Value = -19066.25
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H97").Value = 6350.9
$ws.Range("J97").Value = 6350.9
$ws.Range("L97").Value = 19052.7
$ws.Range("N97").Value = -20044.7
$ws.Range("H102").Value = 40657
$ws.Range("J102").Value = 40657
$ws.Range("L102").Value = 40657
$ws.Range("N102").Value = -47147
$ws.Range("H135").Value = 1219.1212
$ws.Range("I135").Value = 1126.6
$ws.Range("K135").Value = 10139.4
$ws.Range("M135").Value = -7604.4
$ws.Range("H138").Value = 4065.077
$ws.Range("J138").Value = 4870.6665
$ws.Range("L138").Value = 14611.9995
$ws.Range("N138").Value = -24891.9995
$ws.Range("H141").Value = 2599.0952
$ws.Range("I141").Value = 2151.8235
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 6455.470499999999
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = -1275.470499999999
$ws.Range("N141").Value = -23860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1314.5098
$ws.Range("I74").Value = 1355.7858
$ws.Range("J74").Value = 1121.8889
$ws.Range("K74").Value = 1355.7858
$ws.Range("L74").Value = 1121.8889
$ws.Range("M74").Value = -481.7858000000001
$ws.Range("N74").Value = -2869.8889
$ws.Range("H77").Value = 1314.5098
$ws.Range("I77").Value = 1355.7858
$ws.Range("J77").Value = 1121.8889
$ws.Range("K77").Value = 6778.929
$ws.Range("L77").Value = 5609.4445
$ws.Range("M77").Value = -2410.929
$ws.Range("N77").Value = -14345.4445
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3135.7334
$ws.Range("I94").Value = 3215.4546
$ws.Range("J94").Value = 2916.5
$ws.Range("K94").Value = 3215.4546
$ws.Range("L94").Value = 2916.5
$ws.Range("M94").Value = -2764.4546
$ws.Range("N94").Value = -3818.5
$ws.Range("H138").Value = 36510
$ws.Range("J138").Value = 36510
$ws.Range("L138").Value = 36510
$ws.Range("N138").Value = -46790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 940.6667
$ws.Range("I16").Value = 622.5
$ws.Range("K16").Value = 622.5
$ws.Range("M16").Value = -335.5
$ws.Range("H31").Value = 6277.6
$ws.Range("I31").Value = 3004.8
$ws.Range("J31").Value = 7914
$ws.Range("K31").Value = 3004.8
$ws.Range("L31").Value = 7914
$ws.Range("M31").Value = -2709.8
$ws.Range("N31").Value = -8504
$ws.Range("H34").Value = 6277.6
$ws.Range("I34").Value = 3004.8
$ws.Range("J34").Value = 7914
$ws.Range("K34").Value = 3004.8
$ws.Range("L34").Value = 7914
$ws.Range("M34").Value = -2802.8
$ws.Range("N34").Value = -8318
$ws.Range("H113").Value = 940.6667
$ws.Range("I113").Value = 622.5
$ws.Range("K113").Value = 622.5
$ws.Range("M113").Value = 1547.5
$ws.Range("H132").Value = 2098
$ws.Range("I132").Value = 2107.2666
$ws.Range("K132").Value = 6321.7998
$ws.Range("M132").Value = -3791.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 325.14285
$ws.Range("I40").Value = 343.33334
$ws.Range("J40").Value = 292.4
$ws.Range("K40").Value = 1373.33336
$ws.Range("L40").Value = 1169.6
$ws.Range("M40").Value = -1304.33336
$ws.Range("N40").Value = -1307.6
$ws.Range("H140").Value = 3483.75
$ws.Range("I140").Value = 2604.8125
$ws.Range("K140").Value = 7814.4375
$ws.Range("M140").Value = -2634.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1000.75
$ws.Range("I107").Value = 825
$ws.Range("J107").Value = 1176.5
$ws.Range("K107").Value = 825
$ws.Range("L107").Value = 1176.5
$ws.Range("M107").Value = 1095
$ws.Range("N107").Value = -5016.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2355.65
$ws.Range("I22").Value = 2889.818
$ws.Range("K22").Value = 2889.818
$ws.Range("M22").Value = -2594.818
$ws.Range("H27").Value = 2355.65
$ws.Range("I27").Value = 2889.818
$ws.Range("K27").Value = 2889.818
$ws.Range("M27").Value = -2782.818
$ws.Range("H40").Value = 5442.5
$ws.Range("I40").Value = 3968.5652
$ws.Range("K40").Value = 3968.5652
$ws.Range("M40").Value = -3832.5652
$ws.Range("H46").Value = 2667
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 3000.5
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 3000.5
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -3376.5
$ws.Range("H61").Value = 4150.864
$ws.Range("I61").Value = 1574
$ws.Range("K61").Value = 1574
$ws.Range("M61").Value = -1372
$ws.Range("H113").Value = 4150.864
$ws.Range("I113").Value = 1574
$ws.Range("K113").Value = 1574
$ws.Range("M113").Value = 596
$ws.Range("H122").Value = 4873.2573
$ws.Range("I122").Value = 4293.222
$ws.Range("K122").Value = 12879.666
$ws.Range("M122").Value = -10429.666
$ws.Range("H132").Value = 3565.0356
$ws.Range("I132").Value = 3709.2083
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 11127.6249
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -8597.624899999999
$ws.Range("N132").Value = -13160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 50000
$ws.Range("L86").Value = 50000
$ws.Range("N86").Value = -52246
$ws.Range("H89").Value = 50000
$ws.Range("J89").Value = 50000
$ws.Range("L89").Value = 250000
$ws.Range("N89").Value = -261232
$ws.Range("H100").Value = 3898.75
$ws.Range("I100").Value = 3837.6
$ws.Range("K100").Value = 7675.2
$ws.Range("M100").Value = -7134.2
$ws.Range("H138").Value = 81648.664
$ws.Range("J138").Value = 81648.664
$ws.Range("L138").Value = 81648.664
$ws.Range("N138").Value = -91928.664
